# Collapses the three blank "spacer" heading paragraphs that used to sit
# above the EDUCATION heading (part of moving the Education section) into
# a single heading paragraph.
#
# Before:
#   Para A: empty Heading3 paragraph (bookmark "_nigt4dy2guq7")
#   Para B: empty Heading3 paragraph (bookmark "_7q3gk99ea9hw")
#   Para C: empty plain paragraph
#   Para D: Heading3 paragraph containing "EDUCATION" (bookmark "_bxachrh30qxo")
#
# After:
#   Para A: Heading3 paragraph containing "EDUCATION", keeping Para A's own
#           paragraph formatting/bookmark id, but the bookmark is renamed to
#           "_bxachrh30qxo" (Para D's old name) and the "EDUCATION" run
#           (and its bookmark/run content) comes from Para D.
#
# This is what Word does natively when you select from the start of Para A
# through the start of Para D (i.e. every paragraph mark in between,
# including Para A's, B's and C's own marks) and hit Delete: Para A's own
# paragraph mark/formatting survives and Para D's content joins directly
# onto it.

$d = $word.ActiveDocument

$firstPara = $null
$lastPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    # Paragraph.Range.Text includes the trailing paragraph-mark character;
    # strip it before comparing against literal text.
    $text = $p.Range.Text.TrimEnd([char]13)

    if ($null -eq $firstPara) {
        if ($p.Style.NameLocal -eq "Heading 3" -and $text -eq "") {
            $firstPara = $p
        }
    } elseif ($p.Style.NameLocal -eq "Heading 3" -and $text -eq "EDUCATION") {
        $lastPara = $p
        break
    }
}

if ($null -ne $firstPara -and $null -ne $lastPara) {
    # Delete from the start of the first blank heading paragraph through (but
    # not including) the start of the EDUCATION paragraph. This removes the
    # intervening paragraph marks so EDUCATION's own paragraph joins/merges
    # into the first blank heading paragraph, which keeps its own pPr.
    $mergeRange = $d.Range($firstPara.Range.Start, $lastPara.Range.Start)
    $mergeRange.Delete()
}
